$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.621.29'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '2.194.59'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''259.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('D6').Value = '''81.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.77%  '
$ws.Range('D7').Value = '''0.624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.18%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.592'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('D10').Value = '''43.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.86%  '
$ws.Range('D11').Value = '''0.0918'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').Value = '''6.95'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.35%  '
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').Value = '2.522.51'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').Value = '''14.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '2.170.27'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '43.543.00'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').Value = '''0.0000102'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').Value = '''70.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '''5.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('D22').Value = '''2.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +13.85%  '
$ws.Range('D23').Value = '''230.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('E24').Value = '  -4.83%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '''42.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +15.37%  '
$ws.Range('D27').Value = '''10.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +3.15%  '
$ws.Range('D30').Value = '''2.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.13%  '
$ws.Range('D31').Value = '''173.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('D32').Value = '''20.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.34%  '
$ws.Range('D33').Value = '''0.0866'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.24%  '
$ws.Range('E34').Value = '  +4.41%  '
$ws.Range('E35').Value = '  +7.35%  '
$ws.Range('D37').Value = '''4.47'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.03%  '
$ws.Range('E38').Value = '  +5.82%  '
$ws.Range('D39').Value = '''13.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.98%  '
$ws.Range('D40').Value = '''2.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.76%  '
$ws.Range('D41').Value = '''2.09'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('D42').Value = '''63.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.84%  '
$ws.Range('D43').Value = '''5.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.83%  '
$ws.Range('E44').Value = '  +2.17%  '
$ws.Range('D45').Value = '''100.69'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('D46').Value = '''0.0978'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').Value = '''8.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('E48').Value = '  +4.39%  '
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('E50').Value = '  -4.35%  '
$ws.Range('E51').Value = '  +25.52%  '
